# Auto-generated Excel COM-interop script
# Applies the "Updated cryptos list" diff to the worksheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Cell, $Text) {
    # Preserve exact text (including values that look numeric, e.g. "0.5260")
    # by forcing the cell to Text format before assigning, then restoring
    # the original style so no visible formatting changes are introduced.
    $origStyle = $Cell.Style
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = $origStyle
}

$ws.Range("D2").Value = '30.640.96'
$ws.Range("E2").Value = '  +0.56%  '
$ws.Range("D3").Value = '2.114.80'
$ws.Range("E3").Value = '  +0.27%  '
$ws.Range("E4").Value = '  +1.07%  '
Set-TextValue $ws.Range("D5") '348.68'
$ws.Range("E5").Value = '  +3.62%  '
$ws.Range("E6").Value = '  +0.89%  '
Set-TextValue $ws.Range("D7") '0.5260'
$ws.Range("E7").Value = '  +0.39%  '
$ws.Range("E8").Value = '  -1.65%  '
Set-TextValue $ws.Range("D9") '53.71'
$ws.Range("E9").Value = '  +0.74%  '
Set-TextValue $ws.Range("D10") '0.09024'
$ws.Range("E10").Value = '  +1.09%  '
$ws.Range("E11").Value = '  -0.61%  '
Set-TextValue $ws.Range("D12") '24.51'
$ws.Range("E12").Value = '  -0.03%  '
$ws.Range("D13").Value = '2.105.56'
$ws.Range("E13").Value = '  +0.58%  '
Set-TextValue $ws.Range("D14") '6.823'
$ws.Range("E14").Value = '  +0.23%  '
Set-TextValue $ws.Range("D15") '8.034'
$ws.Range("E15").Value = '  +0.42%  '
Set-TextValue $ws.Range("D16") '101.50'
$ws.Range("E16").Value = '  +5.11%  '
$ws.Range("E17").Value = '  +3.23%  '
Set-TextValue $ws.Range("D18") '1.011'
$ws.Range("E18").Value = '  +0.89%  '
Set-TextValue $ws.Range("D19") '0.06709'
$ws.Range("E19").Value = '  +1.18%  '
Set-TextValue $ws.Range("D20") '19.38'
$ws.Range("E20").Value = '  +0.24%  '
Set-TextValue $ws.Range("D21") '1.009'
$ws.Range("E21").Value = '  +0.92%  '
Set-TextValue $ws.Range("D22") '6.304'
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").Value = '30.697.50'
$ws.Range("E23").Value = '  +0.57%  '
Set-TextValue $ws.Range("D24") '12.84'
$ws.Range("E24").Value = '  +3.56%  '
Set-TextValue $ws.Range("D25") '2.386'
$ws.Range("E25").Value = '  +0.96%  '
$ws.Range("D26").Value = '2.361.77'
$ws.Range("E26").Value = '  +1.28%  '
Set-TextValue $ws.Range("D27") '22.35'
$ws.Range("E27").Value = '  -0.19%  '
Set-TextValue $ws.Range("D28") '165.37'
$ws.Range("E28").Value = '  +1.13%  '
Set-TextValue $ws.Range("D29") '2.539'
$ws.Range("E29").Value = '  -1.70%  '
Set-TextValue $ws.Range("D30") '135.84'
$ws.Range("E30").Value = '  +2.20%  '
Set-TextValue $ws.Range("D31") '1.193'
$ws.Range("E31").Value = '  -2.59%  '
Set-TextValue $ws.Range("D32") '0.1075'
$ws.Range("E32").Value = '  +0.02%  '
Set-TextValue $ws.Range("D33") '1.649'
$ws.Range("E33").Value = '  -4.06%  '
Set-TextValue $ws.Range("D34") '6.368'
$ws.Range("E34").Value = '  +2.64%  '
$ws.Range("E35").Value = '  +2.21%  '
Set-TextValue $ws.Range("D36") '5.930'
$ws.Range("E36").Value = '  +6.65%  '
Set-TextValue $ws.Range("D37") '10.24'
$ws.Range("E37").Value = '  -2.65%  '
Set-TextValue $ws.Range("D38") '0.02653'
$ws.Range("E38").Value = '  +2.85%  '
Set-TextValue $ws.Range("D39") '0.06843'
$ws.Range("E39").Value = '  -0.01%  '
Set-TextValue $ws.Range("D40") '0.2314'
$ws.Range("E40").Value = '  +0.35%  '
$ws.Range("E41").Value = '  -2.39%  '
Set-TextValue $ws.Range("D42") '0.6891'
$ws.Range("E42").Value = '  -0.43%  '
Set-TextValue $ws.Range("D43") '1.277'
$ws.Range("E43").Value = '  +2.31%  '
Set-TextValue $ws.Range("D44") '14.75'
$ws.Range("E44").Value = '  +4.85%  '
$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D45") '2.324'
$ws.Range("E45").Value = '  -1.48%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue $ws.Range("D46") '0.6437'
$ws.Range("E46").Value = '  +0.76%  '
Set-TextValue $ws.Range("D47") '3.748'
$ws.Range("E47").Value = '  +2.22%  '
$ws.Range("E48").Value = '  +0.02%  '
$ws.Range("E49").Value = '  +0.09%  '
Set-TextValue $ws.Range("D50") '0.07286'
$ws.Range("E50").Value = '  +2.23%  '
Set-TextValue $ws.Range("D51") '82.29'
$ws.Range("E51").Value = '  -1.86%  '
